$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price / volume values from the symbol-list refresh.
# Cells must stay text (inlineStr) like the source data, so we force
# the Text number format before assigning values that look numeric
# (plain numbers or percentages) -- otherwise Excel auto-converts them
# to real numbers/percent and we lose the original text formatting
# (e.g. trailing zeros, "%" suffix).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "318.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "4.70%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.11"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.158"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.03%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08227"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.96%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.150"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.27%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.010"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.91%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9270"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.00%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "5.43%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1892"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.94%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09157"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.50%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03624"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09918"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.09%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001434"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.64%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005777"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.68%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.462"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.71%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "17.99%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3376"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.48%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.205"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "8.02%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.07%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2191"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.67%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04596"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.53%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001246"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.88%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004732"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.80%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-21.95%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004503"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-5.35%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02005"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "9.44%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04915"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.23%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007787"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.56%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1399"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.07%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "1.03%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002105"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.78%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.11%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006454"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.84%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.11%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.83"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "16.09%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001901"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-5.08%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.11%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.11%"
